# "aggiunti indicatori per semaforo di gruppo"
# Add three new rows (164-166) to the "r AnalysisUnit_Variable" sheet describing
# the new group-level exposure-quota indicators and their related CUSTOMER
# analysis-unit variables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r AnalysisUnit_Variable")

# Data to add: Id/Name (customer indicator) in columns B/C, Variable (indicator) in column F
$newRows = @(
    @{ Row = 164; Id = "CUSTOMER_INDICATOR_QUOTA_EXP_GRP1"; Variable = "INDICATOR_QUOTA_EXP_GRP1" },
    @{ Row = 165; Id = "CUSTOMER_INDICATOR_QUOTA_EXP_GRP2"; Variable = "INDICATOR_QUOTA_EXP_GRP2" },
    @{ Row = 166; Id = "CUSTOMER_INDICATOR_QUOTA_EXP_GRP3"; Variable = "INDICATOR_QUOTA_EXP_GRP3" }
)

# First pass: populate the "Variable" column (F) for all rows so the new
# INDICATOR_QUOTA_EXP_GRP* strings are registered first in the shared-string
# table, then populate the Id/Name columns (B/C) with the CUSTOMER_* strings.
foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("F$r").Value = $item.Variable
    # Highlight the new Variable cell with a yellow fill, as in the source edit
    $ws.Range("F$r").Interior.Color = 65535
}

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = "CREATE/MODIFY"
    $ws.Range("B$r").Value = $item.Id
    $ws.Range("C$r").Value = $item.Id
    $ws.Range("E$r").Value = "CUSTOMER"
}

# Select the last edited cell, matching the author's final cursor position
$ws.Range("B164").Select()
